$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text (not numeric) values are preserved for column D by temporarily
# formatting as text, then resetting the style so no new style is retained.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.056.90'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.00%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.925.92'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.91%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.38%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.62'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.82%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.33%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4600'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.61%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3832'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.83%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07744'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.48%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9838'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.29%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.39'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.20%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.928.76'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.68%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.698'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.62%  '

$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.975'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.07%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07037'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.22%  '

$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '84.35'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.26%  '

$ws.Range("B17").Value = 'BinanceUSD'
$ws.Range("C17").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.002'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.42%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009536'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.38%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.73'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.00%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '29.050.59'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.00%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.344'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.34%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.98'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.75%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.093'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.81%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '158.17'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.73%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '19.14'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.41%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '5.699'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.29%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '118.12'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.72%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.859'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.42%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09342'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.97%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.8676'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.04%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.121'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.05%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.259'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.53%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.049'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.68%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.05713'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.41%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.159'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.64%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.001'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.37%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02052'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.72%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.048'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +13.10%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.547'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.85%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5527'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.77%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1755'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.18%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.000003017'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.09%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '9.407'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.68%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.227'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.59%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5206'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.83%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '11.27'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.18%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06907'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.15%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.783'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.68%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '110.72'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.25%  '

$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '70.09'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.90%  '

